# Fat4-Dchs1 LR-pairs sheet: rebuild the 4x4 (sending cluster x target cluster) block
# with updated NATMI statistics (incl. a new "M2" cluster), per "Natmi following Dr Hou advice".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A..T as in the header row (row 1, untouched).
# Each inner array is one data row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric NATMI columns (E..T). 4 sending clusters x 4 target clusters = 16 rows (rows 2-17).
$rows = @(
    @("ECs", "Fat4", "Dchs1", "ECs", 3, 1, [double]"20.88470466666667", [double]"62.654114", [double]"0.4656575609550316", [double]"0.4656575609550317", 3, 1, [double]"9.616272333333335", [double]"28.848817", [double]"0.3773649226058864", [double]"0.3773649226058863", [double]"200.8330076759043", [double]"1807.497069083138", [double]"0.1757228294506413", [double]"0.1757228294506413"),
    @("ECs", "Fat4", "Dchs1", "FAPs", 3, 1, [double]"20.88470466666667", [double]"62.654114", [double]"0.4656575609550316", [double]"0.4656575609550317", 3, 1, [double]"15.18384366666667", [double]"45.551531", [double]"0.5958493885691961", [double]"0.5958493885691961", [double]"317.1100906831704", [double]"2853.990816148534", [double]"0.2774617729776788", [double]"0.2774617729776788"),
    @("ECs", "Fat4", "Dchs1", "M2", 3, 1, [double]"20.88470466666667", [double]"62.654114", [double]"0.4656575609550316", [double]"0.4656575609550317", 3, 1, [double]"0.2118566666666667", [double]"0.63557", [double]"0.008313749013022724", [double]"0.008313749013022722", [double]"4.424563914997778", [double]"39.82107523498", [double]"0.003871360087796463", [double]"0.003871360087796463"),
    @("ECs", "Fat4", "Dchs1", "sCs", 3, 1, [double]"20.88470466666667", [double]"62.654114", [double]"0.4656575609550316", [double]"0.4656575609550317", 3, 1, [double]"0.4707146666666667", [double]"1.412144", [double]"0.01847193981189477", [double]"0.01847193981189477", [double]"9.830736795601778", [double]"88.476631160416", [double]"0.008601598438915066", [double]"0.008601598438915068"),
    @("FAPs", "Fat4", "Dchs1", "ECs", 3, 1, [double]"21.652497", [double]"64.957491", [double]"0.4827767068068095", [double]"0.4827767068068096", 3, 1, [double]"9.616272333333335", [double]"28.848817", [double]"0.3773649226058864", [double]"0.3773649226058863", [double]"208.216307848683", [double]"1873.946770638147", [double]"0.1821829946000764", [double]"0.1821829946000764"),
    @("FAPs", "Fat4", "Dchs1", "FAPs", 3, 1, [double]"21.652497", [double]"64.957491", [double]"0.4827767068068095", [double]"0.4827767068068096", 3, 1, [double]"15.18384366666667", [double]"45.551531", [double]"0.5958493885691961", [double]"0.5958493885691961", [double]"328.768129440969", [double]"2958.913164968721", [double]"0.2876622055662875", [double]"0.2876622055662876"),
    @("FAPs", "Fat4", "Dchs1", "M2", 3, 1, [double]"21.652497", [double]"64.957491", [double]"0.4827767068068095", [double]"0.4827767068068096", 3, 1, [double]"0.2118566666666667", [double]"0.63557", [double]"0.008313749013022724", [double]"0.008313749013022722", [double]"4.58722583943", [double]"41.28503255487", [double]"0.004013684369725473", [double]"0.004013684369725473"),
    @("FAPs", "Fat4", "Dchs1", "sCs", 3, 1, [double]"21.652497", [double]"64.957491", [double]"0.4827767068068095", [double]"0.4827767068068096", 3, 1, [double]"0.4707146666666667", [double]"1.412144", [double]"0.01847193981189477", [double]"0.01847193981189477", [double]"10.192147907856", [double]"91.72933117070401", [double]"0.008917822270720155", [double]"0.008917822270720157"),
    @("M2", "Fat4", "Dchs1", "ECs", 1, [double]"0.3333333333333333", [double]"0.007427333333333334", [double]"0.022282", [double]"0.0001656041576647308", [double]"0.0001656041576647308", 3, 1, [double]"9.616272333333335", [double]"28.848817", [double]"0.3773649226058864", [double]"0.3773649226058863", [double]"0.07142326004377779", [double]"0.6428093403940001", [double]"6.249320014036413e-05", [double]"6.249320014036412e-05"),
    @("M2", "Fat4", "Dchs1", "FAPs", 1, [double]"0.3333333333333333", [double]"0.007427333333333334", [double]"0.022282", [double]"0.0001656041576647308", [double]"0.0001656041576647308", 3, 1, [double]"15.18384366666667", [double]"45.551531", [double]"0.5958493885691961", [double]"0.5958493885691961", [double]"0.1127754681935556", [double]"1.014979213742", [double]"9.867513608904658e-05", [double]"9.867513608904658e-05"),
    @("M2", "Fat4", "Dchs1", "M2", 1, [double]"0.3333333333333333", [double]"0.007427333333333334", [double]"0.022282", [double]"0.0001656041576647308", [double]"0.0001656041576647308", 3, 1, [double]"0.2118566666666667", [double]"0.63557", [double]"0.008313749013022724", [double]"0.008313749013022722", [double]"0.001573530082222222", [double]"0.01416177074", [double]"1.376791402337615e-06", [double]"1.376791402337615e-06"),
    @("M2", "Fat4", "Dchs1", "sCs", 1, [double]"0.3333333333333333", [double]"0.007427333333333334", [double]"0.022282", [double]"0.0001656041576647308", [double]"0.0001656041576647308", 3, 1, [double]"0.4707146666666667", [double]"1.412144", [double]"0.01847193981189477", [double]"0.01847193981189477", [double]"0.003496154734222223", [double]"0.031465392608", [double]"3.059030032982439e-06", [double]"3.059030032982439e-06"),
    @("sCs", "Fat4", "Dchs1", "ECs", 3, 1, [double]"2.305291666666667", [double]"6.915875", [double]"0.05140012808049411", [double]"0.05140012808049412", 3, 1, [double]"9.616272333333335", [double]"28.848817", [double]"0.3773649226058864", [double]"0.3773649226058863", [double]"22.16831247443056", [double]"199.514812269875", [double]"0.01939660535502831", [double]"0.01939660535502831"),
    @("sCs", "Fat4", "Dchs1", "FAPs", 3, 1, [double]"2.305291666666667", [double]"6.915875", [double]"0.05140012808049411", [double]"0.05140012808049412", 3, 1, [double]"15.18384366666667", [double]"45.551531", [double]"0.5958493885691961", [double]"0.5958493885691961", [double]"35.00318827273611", [double]"315.028694454625", [double]"0.03062673488914078", [double]"0.03062673488914079"),
    @("sCs", "Fat4", "Dchs1", "M2", 3, 1, [double]"2.305291666666667", [double]"6.915875", [double]"0.05140012808049411", [double]"0.05140012808049412", 3, 1, [double]"0.2118566666666667", [double]"0.63557", [double]"0.008313749013022724", [double]"0.008313749013022722", [double]"0.4883914081944444", [double]"4.395522673749999", [double]"0.0004273277640984495", [double]"0.0004273277640984495"),
    @("sCs", "Fat4", "Dchs1", "sCs", 3, 1, [double]"2.305291666666667", [double]"6.915875", [double]"0.05140012808049411", [double]"0.05140012808049412", 3, 1, [double]"0.4707146666666667", [double]"1.412144", [double]"0.01847193981189477", [double]"0.01847193981189477", [double]"1.085134598444444", [double]"9.766211386", [double]"0.0009494600722265697", [double]"0.0009494600722265698")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}